$d = $word.ActiveDocument
$wdParagraph = 4

# The GitHub-description paragraph is immediately followed by a leftover
# empty paragraph (just a blank line) before the "Ly do su dung" heading.
# Locate the GitHub paragraph by its unique text, expand to the full
# paragraph (including its paragraph mark), then grab the paragraph that
# follows it and delete it (since it is blank) to remove the stray blank
# line.
$anchor = $d.Content
$found = $anchor.Find.Execute("GitHub cung cấp nhiều tính năng mạnh mẽ", `
    $true, $false, $false, $false, $false, $true, 1, $false, "", 0)

if (-not $found) {
    throw "Could not find the GitHub description paragraph"
}

[void]$anchor.Expand($wdParagraph)

$blankParagraph = $d.Range($anchor.End, $anchor.End)
[void]$blankParagraph.Expand($wdParagraph)

if ($blankParagraph.Text.Trim().Length -ne 0) {
    throw "Expected the paragraph following the GitHub description to be blank"
}

$blankParagraph.Delete()
